$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add new column J: header + description + values for the new
# "NormalLevelEnemyLocation" / "Data/Level" asset data used for random
# generation of normal level enemy bags.
$ws.Range("J1").Value = "NormalLevelEnemyLocation"
$ws.Range("J2").Value = "Data/Level"
$ws.Range("J3").Value = "Data/Level"
$ws.Range("J4").Value = "Data/Level"
$ws.Range("J5").Value = "Data/Level"

# Move the active selection like in the authored workbook.
$ws.Range("K9").Select()

$wb.Save()
